# Updated notebook, reran simulation
#
# The underlying "label" list used to build this sheet gained two new
# entries ("Holden" and "Rizzie Spiral") inserted near its front, and the
# simulation now also reports two additional rows (28 and 29) that weren't
# present before. We rebuild the row-label column (B) and the plane-header
# row (row 2, C:W) from the updated master label list, and append the two
# new data rows (30/31) with the same per-reflection value of 1 used
# throughout the grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated master label list (after the two new entries were inserted).
$labels = @(
    "HKL",
    "Spiral5",
    "Holden",
    "Rizzie Spiral",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Matthies Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex",
    "[2, 0, 0]",
    "[2, 2, 0]",
    "[3, 3, 3]",
    "[4, 2, 0]",
    "[4, 0, 0]",
    "[4, 2, 2]",
    "[5, 1, 1]",
    "[1, 1, 1]",
    "[2, 2, 2]",
    "[3, 3, 1]",
    "[3, 1, 1]",
    "1Pair-A",
    "1Pair-B",
    "2Pairs-A",
    "2Pairs-B",
    "3Pairs-A",
    "3Pairs-B",
    "3Pairs-C",
    "4Pairs",
    "5A4F",
    "MaxUnique"
)

# --- First: append the two new simulated rows (30 & 31), copying the
# formatting of the last existing row so the bold/border/center style
# carries over exactly like the rest of column A. ---
$ws.Range("A29:W29").Copy($ws.Range("A30:W30"))
$ws.Range("A29:W29").Copy($ws.Range("A31:W31"))

# --- Row labels: column A is the zero-based row index, column B is the
# corresponding entry from the (now longer) label list. Rows 2-31 map to
# label indices 0-29. ---
for ($r = 2; $r -le 31; $r++) {
    $i = $r - 2
    $ws.Cells.Item($r, 1).Value() = $i
    $ws.Cells.Item($r, 2).Value() = $labels[$i]
}

# --- Header row 2 (columns C:W) are the remaining label-list entries
# (indices 30-50), one per "plane" column. ---
for ($c = 3; $c -le 23; $c++) {
    $i = $c + 27
    $ws.Cells.Item(2, $c).Value() = $labels[$i]
}

# --- Data grid for the two new rows: every reflection column is 1, same
# as every other row in the sheet. ---
for ($c = 3; $c -le 23; $c++) {
    $ws.Cells.Item(30, $c).Value() = 1
    $ws.Cells.Item(31, $c).Value() = 1
}
